$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text values (e.g. "28.299.26",
# "  +1.22%  ") rather than numbers. Force the cells to Text format first so Excel
# does not reinterpret numeric-looking strings (like "1.006") as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.299.26'
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '1.886.83'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = '314.79'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '0.5145'
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").Value = '0.3924'
$ws.Range("E8").Value = '  +3.17%  '
$ws.Range("D9").Value = '0.08378'
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").Value = '41.74'
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("D12").Value = '6.241'
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("D13").Value = '20.77'
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("D14").Value = '1.881.32'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '7.303'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").Value = '1.007'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '91.63'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").Value = '0.06679'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D22").Value = '6.064'
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("D23").Value = '28.349.57'
$ws.Range("E23").Value = '  +1.16%  '
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '2.284'
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.093.56'
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.531'
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '159.42'
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '20.70'
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '125.87'
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.1069'
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '1.051'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.912'
$ws.Range("E33").Value = '  +5.85%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '3.605'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").Value = '9.773'
$ws.Range("E35").Value = '  +2.25%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02469'
$ws.Range("E36").Value = '  +2.67%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '0.06596'
$ws.Range("E37").Value = '  +1.34%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = '0.2201'
$ws.Range("E38").Value = '  +2.34%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = '1.217'
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.6549'
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").Value = '5.033'
$ws.Range("E41").Value = '  +3.67%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.232'
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '11.34'
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = '0.6180'
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.16'
$ws.Range("E45").Value = '  +1.40%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.290'
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.689'
$ws.Range("E47").Value = '  +0.91%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.023'
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.240'
$ws.Range("E49").Value = '  +2.77%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '121.74'
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '79.21'
$ws.Range("E51").Value = '  -0.37%  '

$wb.Save()
